$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","H","K","L","N")

$data = @(
    @(19.01164980593784, 5.266203629647237, 8.475935357025193, 10.66851330593078, 42.06301780364799, 7.344005520526261, 15.21615303703812, 10.25156721744367, 23.72442181002088),
    @(18.82691776739079, 5.115554408922794, 8.487262257346668, 10.67603973725415, 41.90554509757459, 7.344005520526261, 15.09156237742939, 10.24509325233009, 23.75357729588041),
    @(18.7179461065581, 5.019570367264432, 8.49442753970146, 10.68232289721508, 41.81724240539624, 7.344005520526261, 15.01872627120869, 10.24307119016048, 23.77330635176207),
    @(18.67470715159067, 4.979606918859592, 8.497400636014982, 10.68530135535586, 41.78338688379616, 7.344005520526261, 14.98999542398287, 10.24273923136116, 23.78180538484989),
    @(18.66759923966785, 4.972920571241506, 8.497897535942506, 10.68582117493794, 41.77789428633522, 7.344005520526261, 14.98528291467458, 10.2427138509357, 23.78324437496583),
    @(18.71735818097925, 5.019034806837348, 8.494467420262662, 10.68236137311064, 41.81677717554207, 7.344005520526261, 15.01833491124054, 10.24306471990621, 23.77341911347399),
    @(18.94705903817247, 5.214998555492127, 8.479797330255094, 10.67076359213363, 42.00699327286463, 7.344005520526261, 15.17245207988038, 10.24893030439076, 23.7340951807831),
    @(19.4306372112221, 5.570566680476789, 8.452687818551992, 10.66119867400916, 42.44556865764091, 7.344005520526261, 15.50236473804663, 10.2758764525803, 23.67149296911021),
    @(19.80303359183134, 5.81303747109269, 8.433764101216113, 10.66218998278867, 42.80629427369976, 7.344005520526261, 15.7596832402028, 10.30500359385088, 23.63436293133038),
    @(19.97547845760167, 5.919047021755785, 8.425367191686352, 10.66437671337544, 42.97839853697947, 7.344005520526261, 15.8795468328966, 10.3202562980099, 23.61939988494012),
    @(20.04115772976651, 5.958556980262408, 8.422217649610159, 10.66545368373633, 43.0446876430928, 7.344005520526261, 15.92530176459486, 10.32631755224988, 23.6140112532013),
    @(20.02699664530562, 5.950076223150838, 8.422894621163737, 10.66521068130585, 43.03036199338282, 7.344005520526261, 15.91543199915291, 10.32499950427864, 23.61515944164048),
    @(19.9808746904666, 5.922310314703044, 8.425107473381122, 10.66446033242712, 42.98383000154216, 7.344005520526261, 15.88330402254088, 10.32074925730532, 23.6189509944061),
    @(19.95267115444551, 5.905219930628989, 8.426466833153224, 10.66403311375273, 42.95547221205208, 7.344005520526261, 15.86367111359954, 10.318182942393, 23.62130958573551),
    @(19.79182004576072, 5.80602169174831, 8.434317096071453, 10.66208194603288, 42.79520508890733, 7.344005520526261, 15.75190309001426, 10.30404687682316, 23.63537962091703),
    @(19.69387978337037, 5.744055304188598, 8.439186998829646, 10.66132908270618, 42.69891353088764, 7.344005520526261, 15.68402895337468, 10.2958860261675, 23.64450509514237),
    @(19.63783681698181, 5.708010785781524, 8.442007967760306, 10.6610594849348, 42.64428516450987, 7.344005520526261, 15.64525621060856, 10.29138078995766, 23.64993523949881),
    @(19.61891317897813, 5.695738031514388, 8.442966529627506, 10.66099629224534, 42.62591980289406, 7.344005520526261, 15.63217539394229, 10.28988787272695, 23.65180493971757),
    @(19.70427613742752, 5.750693575230981, 8.438666528497357, 10.66139231668066, 42.70908593429404, 7.344005520526261, 15.69122695040253, 10.29673525397442, 23.6435148947505),
    @(19.99441203442764, 5.930483151286945, 8.424456687734788, 10.6646739791038, 42.99746754738029, 7.344005520526261, 15.89273118694932, 10.32198993459263, 23.61782978797187),
    @(20.18620938580685, 6.044286846275609, 8.415345567364316, 10.66826912210876, 43.19243466392796, 7.344005520526261, 16.02653627494011, 10.34015737936412, 23.60266101968242),
    @(20.0836638472228, 5.98389109788751, 8.420192334434407, 10.66621787855136, 43.0877948714553, 7.344005520526261, 15.95494174347685, 10.33030992193643, 23.61060871718649),
    @(19.6995751183359, 5.747693717433744, 8.438901767075027, 10.66136322010884, 42.70448471127543, 7.344005520526261, 15.68797195740217, 10.29635073683333, 23.64396199185632),
    @(19.29658223897301, 5.477589029895853, 8.45984590930092, 10.66237652940642, 42.32004519311515, 7.344005520526261, 15.41034627501635, 10.26694110400127, 23.68687315864325),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}
